$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    2 = @{ B = "1111111111011011100010010000000100001101100011011101101111111111"; C = 0.015064 }
    3 = @{ B = "1111111111011011100010010000000100001101100011111101101111111111"; C = 0.031337; D = 1 }
    4 = @{ B = "1111111111011011100010010000000100001101100011011101101111111111"; C = 0 }
    5 = @{ B = "1111111111011011100010010000000100001101100011011101101111111111"; C = 0.029687 }
    6 = @{ B = "1111111111011011100010010000000100001101100011111101101111111111"; C = 0.015019; D = 1 }
    7 = @{ B = "1111111111011011100010010000000100001101100011111101101111111111"; D = 1 }
    8 = @{ B = "1111111111011011100010010000000100001101100011111101101111111111"; D = 1 }
    9 = @{ B = "1111111111011011100010010000000100001101100011111101101111111111"; D = 1 }
    10 = @{ B = "1111111111011011100010010000000100001101100011111101101111111111"; D = 1 }
    11 = @{ B = "1111111111011011100010010000000100001101100011111101101111111111"; D = 1 }
    12 = @{ B = "1111111111011011100010010000000100001101100011111101101111111111"; D = 1 }
    13 = @{ B = "1111111111011011100010010000000100001101100011111101101111111111"; D = 1 }
    14 = @{ B = "1111111111011011100010010000000100001101100011111101101111111111"; D = 1 }
    15 = @{ B = "1111111111011011100010010000000100001101100011111101101111111111"; D = 1 }
    16 = @{ B = "1111111111011011100010010000000100001101100011111101101111111111"; D = 1 }
    17 = @{ B = "1111111111011011100010010000000100001101100011111101101111111111"; D = 1 }
    18 = @{ B = "1111111111011011100010010000000100001101100011111101101111111111"; C = 0; D = 1 }
    19 = @{ B = "1111111111011011100010010000000100001101100011111101101111111111"; D = 1 }
    20 = @{ B = "1111111111011011100010010000000100001101100011111101101111111111"; D = 1 }
    21 = @{ B = "1111111111011011100010010000000100001101100011111101101111111111"; D = 1 }
    22 = @{ B = "1111111111011011100010010000000100001101100011111101101111111111"; D = 1 }
    23 = @{ B = "1111111111011011100010010000000100001101100011111101101111111111"; D = 1 }
    24 = @{ B = "1111111111011011100010010000000100001101100011111101101111111111"; C = 0.0157; D = 1 }
    25 = @{ B = "1111111111011011100010010000000100001101100011111101101111111111"; D = 1 }
    26 = @{ B = "1111111111011011100010010000000100001101100011111101101111111111"; D = 1 }
    27 = @{ B = "1111111111011011100010010000000100001101100011011101101111111111" }
    28 = @{ B = "1111111111011011100010010000000100001101100011011101101111111111" }
    29 = @{ B = "1111111111011011100010010000000100001101100011011101101111111111" }
    30 = @{ B = "1111111111011011100010010000000100001101100011011101101111111111"; C = 0 }
    31 = @{ B = "1111111111011011100010010000000100001101100011011101101111111111" }
    32 = @{ B = "1111111111011011100010010000000100001101100011011101101111111111" }
    33 = @{ B = "1111111111011011100010010000000100001101100011011101101111111111" }
    34 = @{ B = "1111111111011011100010010000000100001101100011011101101111111111" }
    35 = @{ B = "1111111111011011100010010000000100001101100011011101101111111111" }
    36 = @{ B = "1111111111011011100010010000000100001101100011011101101111111111" }
    37 = @{ B = "1111111111011011100010010000000100001101100011011101101111111111" }
    38 = @{ B = "1111111111011011100010010000000100001101100011011101101111111111"; C = 0.015626 }
    39 = @{ B = "1111111111011011100010010000000100001101100011011101101111111111" }
    40 = @{ B = "1111111111011011100010010000000100001101100011011101101111111111" }
    41 = @{ B = "1111111111011011100010010000000100001101100011011101101111111111" }
    42 = @{ B = "1111111111011011100010010000000100001101100011011101101111111111" }
    43 = @{ B = "1111111111011011100010010000000100001101100011011101101111111111"; C = 0 }
    44 = @{ B = "1111111111011011100010010000000100001101100011011101101111111111" }
    45 = @{ B = "1111111111011011100010010000000100001101100011011101101111111111" }
    46 = @{ B = "1111111111011011100010010000000100001101100011011101101111111111" }
    47 = @{ B = "1111111111011011100010010000000100001101100011011101101111111111" }
    48 = @{ B = "1111101111011011100010010000000100001101100011011101101111111111"; C = 0.015629; D = 1 }
    49 = @{ B = "1111111111011011100010010000000100001101100011011101101111111111"; C = 0.015635 }
    50 = @{ B = "1111101111011011100010010000000100001101100011011101101111111111"; C = 0.015637; D = 1 }
    51 = @{ B = "1111111111011011100010010000000100001101100011011101101111111111"; C = 0.015609 }
    52 = @{ B = "1111111111011011100010010000000100001101100011011101101111111111"; C = 0 }
    53 = @{ B = "1111111111011011100010010000000100001101100011011101101111111111"; C = 0.015625 }
    54 = @{ B = "1111111111011011100010010000000100001101100011011101101111111111"; C = 0.015635 }
    55 = @{ B = "1111111111011011100010010000000100001101100011011101101111111111"; C = 0.015616 }
    56 = @{ B = "1111111111011011100010010000000100001101100011111101111111111111"; C = 0.015631; D = 2 }
    57 = @{ B = "1111111111011011100011110000001100001101100011111101111111111111"; C = 0.015616; D = 5 }
    58 = @{ B = "1111111111011111100011110000101100001111100011111101111111111111"; C = 0.016457; D = 8 }
    59 = @{ B = "1111111111011111100011110000111100001111100011111101111111111111"; C = 0.014163; D = 9 }
    60 = @{ B = "1111111111011111100011110000111100001111100011111101111111111111"; C = 0.012019; D = 9 }
    61 = @{ B = "1111111111011111100011110000111100101111100011111101111111111111"; C = 0; D = 10 }
    62 = @{ B = "1111111111011111100011110000111100101111000011111101111111111111"; C = 0.015641; D = 11 }
    63 = @{ B = "1111111111011111100011110000111100101111000011111101111111111111"; C = 0.027684; D = 11 }
    64 = @{ B = "1111111111011111000011110000111100101111000011111101111111111111"; C = 0.005097; D = 12 }
}

foreach ($r in $changes.Keys) {
    $row = $changes[$r]
    $rowNum = [int]$r
    if ($row.ContainsKey("B")) {
        $ws.Cells.Item($rowNum, 2).NumberFormat = "@"
        $ws.Cells.Item($rowNum, 2).Value = $row.B
    }
    if ($row.ContainsKey("C")) {
        $ws.Cells.Item($rowNum, 3).Value = $row.C
    }
    if ($row.ContainsKey("D")) {
        $ws.Cells.Item($rowNum, 4).Value = $row.D
    }
}

Write-Host "Applied $($changes.Count) row updates"